$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.174785100286533
$ws.Range("C2").Value = 0.5959885386819485
$ws.Range("J2").Value = 0.0171919770773639
$ws.Range("P2").Value = 0.1404011461318052
$ws.Range("S2").Value = 0.07163323782234957
$ws.Range("B3").Value = 0.004761904761904762
$ws.Range("C3").Value = 0.004761904761904762
$ws.Range("J3").Value = 0.03333333333333333
$ws.Range("P3").Value = 0.8142857142857143
$ws.Range("S3").Value = 0.1428571428571428
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.7833333333333333
$ws.Range("S4").Value = 0.1333333333333333
$ws.Range("B6").Value = 0.05583756345177665
$ws.Range("D6").Value = 0.03045685279187817
$ws.Range("E6").Value = 0.005076142131979695
$ws.Range("F6").Value = 0.05076142131979695
$ws.Range("J6").Value = 0.2030456852791878
$ws.Range("O6").Value = 0.03045685279187817
$ws.Range("Q6").Value = 0.2487309644670051
$ws.Range("R6").Value = 0.06598984771573604
$ws.Range("S6").Value = 0.3096446700507614
$ws.Range("B7").Value = 0.1121951219512195
$ws.Range("D7").Value = 0.03414634146341464
$ws.Range("E7").Value = 0.004878048780487805
$ws.Range("F7").Value = 0.06829268292682927
$ws.Range("J7").Value = 0.1512195121951219
$ws.Range("O7").Value = 0.01463414634146342
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.08292682926829269
$ws.Range("S7").Value = 0.3317073170731707
$ws.Range("B8").Value = 0.1084598698481562
$ws.Range("D8").Value = 0.03036876355748373
$ws.Range("E8").Value = 0.002169197396963124
$ws.Range("F8").Value = 0.03904555314533623
$ws.Range("J8").Value = 0.1127982646420824
$ws.Range("O8").Value = 0.02386117136659436
$ws.Range("Q8").Value = 0.2516268980477224
$ws.Range("R8").Value = 0.0455531453362256
$ws.Range("S8").Value = 0.386117136659436
$ws.Range("B9").Value = 0.09420289855072464
$ws.Range("D9").Value = 0.02536231884057971
$ws.Range("E9").Value = 0.003623188405797101
$ws.Range("F9").Value = 0.02536231884057971
$ws.Range("J9").Value = 0.1413043478260869
$ws.Range("O9").Value = 0.02173913043478261
$ws.Range("Q9").Value = 0.2282608695652174
$ws.Range("R9").Value = 0.05072463768115942
$ws.Range("S9").Value = 0.4094202898550725
$ws.Range("B10").Value = 0.1328358208955224
$ws.Range("D10").Value = 0.02164179104477612
$ws.Range("E10").Value = 0.002238805970149254
$ws.Range("F10").Value = 0.07388059701492537
$ws.Range("J10").Value = 0.1111940298507463
$ws.Range("O10").Value = 0.006716417910447761
$ws.Range("Q10").Value = 0.2223880597014925
$ws.Range("R10").Value = 0.07164179104477612
$ws.Range("S10").Value = 0.3574626865671642
$ws.Range("G11").Value = 0.1186943620178042
$ws.Range("J11").Value = 0.09495548961424333
$ws.Range("K11").Value = 0.1810089020771513
$ws.Range("L11").Value = 0.5816023738872403
$ws.Range("S11").Value = 0.02373887240356083
$ws.Range("G12").Value = 0.68
$ws.Range("J12").Value = 0.26
$ws.Range("K12").Value = 0.005
$ws.Range("L12").Value = 0.025
$ws.Range("S12").Value = 0.03
$ws.Range("G13").Value = 0.62
$ws.Range("J13").Value = 0.34
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.004950495049504951
$ws.Range("H15").Value = 0.1138613861386139
$ws.Range("I15").Value = 0.08415841584158416
$ws.Range("J15").Value = 0.3663366336633663
$ws.Range("K15").Value = 0.103960396039604
$ws.Range("M15").Value = 0.01485148514851485
$ws.Range("O15").Value = 0.0396039603960396
$ws.Range("S15").Value = 0.2722772277227723
$ws.Range("F16").Value = 0.007326007326007326
$ws.Range("H16").Value = 0.1978021978021978
$ws.Range("I16").Value = 0.0989010989010989
$ws.Range("J16").Value = 0.3296703296703297
$ws.Range("K16").Value = 0.1245421245421245
$ws.Range("M16").Value = 0.02197802197802198
$ws.Range("N16").Value = 0.007326007326007326
$ws.Range("O16").Value = 0.05494505494505494
$ws.Range("S16").Value = 0.1575091575091575
$ws.Range("F17").Value = 0.0071301247771836
$ws.Range("H17").Value = 0.1443850267379679
$ws.Range("I17").Value = 0.1087344028520499
$ws.Range("J17").Value = 0.427807486631016
$ws.Range("K17").Value = 0.1140819964349376
$ws.Range("M17").Value = 0.0196078431372549
$ws.Range("O17").Value = 0.0392156862745098
$ws.Range("S17").Value = 0.1390374331550802
$ws.Range("F18").Value = 0.0124223602484472
$ws.Range("H18").Value = 0.1801242236024845
$ws.Range("I18").Value = 0.124223602484472
$ws.Range("J18").Value = 0.3975155279503105
$ws.Range("K18").Value = 0.08074534161490683
$ws.Range("M18").Value = 0.01863354037267081
$ws.Range("O18").Value = 0.06832298136645963
$ws.Range("S18").Value = 0.1180124223602484
$ws.Range("F19").Value = 0.01196709050112191
$ws.Range("H19").Value = 0.206432311144353
$ws.Range("I19").Value = 0.112191473448018
$ws.Range("J19").Value = 0.3455497382198953
$ws.Range("K19").Value = 0.1032161555721765
$ws.Range("M19").Value = 0.02318623784592371
$ws.Range("N19").Value = 0.001495886312640239
$ws.Range("O19").Value = 0.05908750934928945
$ws.Range("S19").Value = 0.1368735976065819
